$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 649
$ws.Range("C2").Value = 368
$ws.Range("B3").Value = 999
$ws.Range("C3").Value = 1993
$ws.Range("B4").Value = 73
$ws.Range("C4").Value = 37
$ws.Range("C5").Value = 110
